# aggiornamento 15, 16, 17 marzo
# Add three new data rows (227-229) below the existing data, continuing the
# daily series with the same formatting as the last existing row (226).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy the formatting (number format, font, borders, alignment) of the last
# existing row so the new rows look identical to the rest of the table,
# while reusing the existing style definitions instead of creating new ones.
$ws.Range("A226:D226").Copy()
$ws.Range("A227:D229").PasteSpecial(-4122)  # xlPasteFormats

# Row 227: 2021-04-15 (serial 44301)
$ws.Range("A227").Value = 44301
$ws.Range("B227").Value = 0
$ws.Range("C227").Value = 1
$ws.Range("D227").Value = 109.1703056768559

# Row 228: 2021-04-16 (serial 44302)
$ws.Range("A228").Value = 44302
$ws.Range("B228").Value = 1
$ws.Range("C228").Value = 1
$ws.Range("D228").Value = 109.1703056768559

# Row 229: 2021-04-17 (serial 44303)
$ws.Range("A229").Value = 44303
$ws.Range("B229").Value = 0
$ws.Range("C229").Value = 1
$ws.Range("D229").Value = 109.1703056768559
